# This script applies the edit described by the diff:
#  1. Near the end of the document, remove the paragraph that duplicates the
#     bold title ("Play Balloonies for Free: Review and Game Overview | IGT")
#     and change the remaining (formerly meta-description, italic) paragraph's
#     text to the new image-generation "Prompt: ..." text.
#  2. Insert a new paragraph right after the document's first paragraph
#     (the Heading1 title) containing a bold "Meta description" run followed
#     by a normal run with the old meta-description text.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Step 1: handle the two paragraphs near the end of the document.
# ---------------------------------------------------------------------------

$titleText = "Play Balloonies for Free: Review and Game Overview | IGT"

# Locate the SECOND occurrence of the title text (the first is the Heading1
# at the very start of the document).
$search = $d.Content
$search.Find.Execute($titleText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$search.Collapse(0)
$search.Find.Execute($titleText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null

# $search now covers just the matched text; expand it to the whole paragraph
# (the bold duplicate-title paragraph) and delete the paragraph entirely,
# including its paragraph mark.
$dupTitlePara = $search.Paragraphs(1)
$dupTitlePara.Range.Delete()

# The paragraph that follows (previously holding the italic meta-description
# text) is now the last paragraph in the document. Replace its text (but not
# its paragraph mark, so the italic run formatting is preserved) with the new
# prompt text.
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$promptRange = $lastPara.Range
$promptRange.MoveEnd(1, -1) | Out-Null
$promptRange.Text = "Prompt: Create a feature image in cartoon style for the game ""Balloonies"" that features a happy Maya warrior with glasses. The Maya warrior should be floating amongst the colorful animal-shaped balloons that are the main theme of the game. Make sure to incorporate elements from the game such as the star-shaped balloon, red bonus balloon, and the balloon-shaped icon featuring a Hedgehog that functions as the Wild symbol. The image should be bright and vibrant, catching the attention of potential players and highlighting the playful and entertaining nature of the game."

# ---------------------------------------------------------------------------
# Step 2: insert the new "Meta description" paragraph after the first
# paragraph (the document title, styled Heading1).
# ---------------------------------------------------------------------------

$firstPara = $d.Paragraphs(1)
$insertPoint = $firstPara.Range
$insertPoint.Collapse(0)
$insertPoint.InsertParagraphAfter()

$metaPara = $d.Paragraphs(2)
# The new paragraph inherits the Heading1 style from the paragraph it split
# off from; switch it back to the Normal body style (matching Heading1's
# "style for following paragraph").
$metaPara.Style = "Normal"

$metaRange = $metaPara.Range
$metaRange.MoveEnd(1, -1) | Out-Null

# Leading empty run (matches the empty "<w:r/>" seen throughout the document).
$metaRange.InsertAfter("")

# Bold "Meta description" run.
$metaRange.Collapse(0)
$metaRange.InsertAfter("Meta description")
$metaRange.Bold = 1

# Trailing, non-bold run with the rest of the sentence.
$metaRange.Collapse(0)
$metaRange.InsertAfter(": Learn about the exciting bonus features and unique gameplay of Balloonies by IGT, and play it for free before trying your luck with real money.")
$metaRange.Bold = 0
